# Scheduled runner update: refresh currentAveragePrice / Leve price / profit
# columns (H:N) for a batch of leves across the ALC, ARM, BSM, CRP, CUL, GSM,
# LTW and WVR sheets with newly pulled market-board data.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H118").Value = 3733.7
$ws.Range("I118").Value = 440
$ws.Range("J118").Value = 4831.6
$ws.Range("K118").Value = 1320
$ws.Range("L118").Value = 14494.8
$ws.Range("M118").Value = 337
$ws.Range("N118").Value = -17808.8

$ws.Range("H135").Value = 18519012
$ws.Range("I135").Value = 402.82352
$ws.Range("J135").Value = 50000650
$ws.Range("K135").Value = 3625.41168
$ws.Range("L135").Value = 450005850
$ws.Range("M135").Value = -1090.41168
$ws.Range("N135").Value = -450010920

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 17144.436
$ws.Range("I2").Value = 19284
$ws.Range("J2").Value = 6981.5
$ws.Range("K2").Value = 19284
$ws.Range("L2").Value = 6981.5
$ws.Range("M2").Value = -19171
$ws.Range("N2").Value = -7207.5

$ws.Range("H45").Value = 1505.8182
$ws.Range("I45").Value = 1275
$ws.Range("J45").Value = 1637.7142
$ws.Range("K45").Value = 1275
$ws.Range("L45").Value = 1637.7142
$ws.Range("M45").Value = -898
$ws.Range("N45").Value = -2391.7142

$ws.Range("H61").Value = 8109850.5
$ws.Range("I61").Value = 9681892
$ws.Range("J61").Value = 6537809.5
$ws.Range("K61").Value = 9681892
$ws.Range("L61").Value = 6537809.5
$ws.Range("M61").Value = -9681680
$ws.Range("N61").Value = -6538233.5

$ws.Range("H63").Value = 2000
$ws.Range("I63").Value = 1833.3334
$ws.Range("J63").Value = 2166.6667
$ws.Range("K63").Value = 1833.3334
$ws.Range("L63").Value = 2166.6667
$ws.Range("M63").Value = -1147.3334
$ws.Range("N63").Value = -3538.6667

$ws.Range("H66").Value = 2000
$ws.Range("I66").Value = 1833.3334
$ws.Range("J66").Value = 2166.6667
$ws.Range("K66").Value = 9166.666999999999
$ws.Range("L66").Value = 10833.3335
$ws.Range("M66").Value = -5734.666999999999
$ws.Range("N66").Value = -17697.3335

$ws.Range("H116").Value = 17144.436
$ws.Range("I116").Value = 19284
$ws.Range("J116").Value = 6981.5
$ws.Range("K116").Value = 19284
$ws.Range("L116").Value = 6981.5
$ws.Range("M116").Value = -16990
$ws.Range("N116").Value = -11569.5

$ws.Range("H132").Value = 15420239
$ws.Range("I132").Value = 13378838
$ws.Range("J132").Value = 19660070
$ws.Range("K132").Value = 40136514
$ws.Range("L132").Value = 58980210
$ws.Range("M132").Value = -40133984
$ws.Range("N132").Value = -58985270

$ws.Range("H136").Value = 8109850.5
$ws.Range("I136").Value = 9681892
$ws.Range("J136").Value = 6537809.5
$ws.Range("K136").Value = 29045676
$ws.Range("L136").Value = 19613428.5
$ws.Range("M136").Value = -29043126
$ws.Range("N136").Value = -19618528.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 17144.436
$ws.Range("I3").Value = 19284
$ws.Range("J3").Value = 6981.5
$ws.Range("K3").Value = 19284
$ws.Range("L3").Value = 6981.5
$ws.Range("M3").Value = -19170
$ws.Range("N3").Value = -7209.5

$ws.Range("H107").Value = 575
$ws.Range("I107").Value = 566.6667
$ws.Range("J107").Value = 600
$ws.Range("K107").Value = 566.6667
$ws.Range("L107").Value = 600
$ws.Range("M107").Value = 1353.3333
$ws.Range("N107").Value = -4440

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 10000001
$ws.Range("J6").Value = 0
$ws.Range("L6").Value = 0
$ws.Range("N6").ClearContents()

$ws.Range("H17").Value = 2980
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 2980
$ws.Range("K17").Value = 0
$ws.Range("M17").Value = 2980
$ws.Range("N17").Value = -3328
$ws.Range("L17").ClearContents()

$ws.Range("H50").Value = 18092
$ws.Range("J50").Value = 18092
$ws.Range("L50").Value = 18092
$ws.Range("N50").Value = -19342

$ws.Range("H58").Value = 1003.86664
$ws.Range("I58").Value = 906.1739
$ws.Range("J58").Value = 1324.8572
$ws.Range("K58").Value = 906.1739
$ws.Range("L58").Value = 1324.8572
$ws.Range("M58").Value = -703.1739
$ws.Range("N58").Value = -1730.8572

$ws.Range("H74").Value = 29484.5
$ws.Range("J74").Value = 29484.5
$ws.Range("L74").Value = 29484.5
$ws.Range("N74").Value = -31232.5

$ws.Range("H77").Value = 29484.5
$ws.Range("J77").Value = 29484.5
$ws.Range("L77").Value = 88453.5
$ws.Range("N77").Value = -97189.5

$ws.Range("H132").Value = 2176560.2
$ws.Range("I132").Value = 3334671.5
$ws.Range("J132").Value = 5101.5
$ws.Range("K132").Value = 10004014.5
$ws.Range("L132").Value = 15304.5
$ws.Range("M132").Value = -10001484.5
$ws.Range("N132").Value = -20364.5

$ws.Range("H134").Value = 1257441.5
$ws.Range("I134").Value = 6532.316
$ws.Range("J134").Value = 3085693.5
$ws.Range("K134").Value = 19596.948
$ws.Range("L134").Value = 9257080.5
$ws.Range("M134").Value = -17061.948
$ws.Range("N134").Value = -9262150.5

$ws.Range("H136").Value = 1003.86664
$ws.Range("I136").Value = 906.1739
$ws.Range("J136").Value = 1324.8572
$ws.Range("K136").Value = 2718.5217
$ws.Range("L136").Value = 3974.5716
$ws.Range("M136").Value = -168.5217000000002
$ws.Range("N136").Value = -9074.571599999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H32").Value = 1224
$ws.Range("I32").Value = 1364.6666
$ws.Range("J32").Value = 380
$ws.Range("K32").Value = 4093.9998
$ws.Range("L32").Value = 1140
$ws.Range("M32").Value = -3810.9998
$ws.Range("N32").Value = -1706

$ws.Range("H33").Value = 266.66666
$ws.Range("I33").Value = 0
$ws.Range("J33").Value = 266.66666
$ws.Range("K33").Value = 0
$ws.Range("M33").Value = 1599.99996
$ws.Range("N33").Value = -2165.99996
$ws.Range("L33").ClearContents()

$ws.Range("H34").Value = 1841.65
$ws.Range("I34").Value = 861.4286
$ws.Range("J34").Value = 2369.4614
$ws.Range("K34").Value = 2584.2858
$ws.Range("L34").Value = 7108.3842
$ws.Range("M34").Value = -2500.2858
$ws.Range("N34").Value = -7276.3842

$ws.Range("H38").Value = 14341.286
$ws.Range("I38").Value = 16726.5
$ws.Range("J38").Value = 30
$ws.Range("K38").Value = 50179.5
$ws.Range("L38").Value = 90
$ws.Range("M38").Value = -49832.5
$ws.Range("N38").Value = -784

$ws.Range("H39").Value = 2672.2222
$ws.Range("J39").Value = 2950
$ws.Range("L39").Value = 8850
$ws.Range("N39").Value = -9438

$ws.Range("H40").Value = 4616.5864
$ws.Range("I40").Value = 9561.909
$ws.Range("J40").Value = 1594.4445
$ws.Range("K40").Value = 38247.636
$ws.Range("L40").Value = 6377.778
$ws.Range("M40").Value = -38178.636
$ws.Range("N40").Value = -6515.778

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 15359.143
$ws.Range("I126").Value = 20840
$ws.Range("J126").Value = 1657
$ws.Range("K126").Value = 62520
$ws.Range("L126").Value = 4971
$ws.Range("M126").Value = -60050
$ws.Range("N126").Value = -9911

$ws.Range("H132").Value = 24411148
$ws.Range("I132").Value = 23001012
$ws.Range("J132").Value = 25977966
$ws.Range("K132").Value = 69003036
$ws.Range("L132").Value = 77933898
$ws.Range("M132").Value = -69000506
$ws.Range("N132").Value = -77938958

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1611.2778
$ws.Range("I7").Value = 1437.125
$ws.Range("J7").Value = 3004.5
$ws.Range("K7").Value = 1437.125
$ws.Range("L7").Value = 3004.5
$ws.Range("M7").Value = -1325.125
$ws.Range("N7").Value = -3228.5

$ws.Range("H22").Value = 27785296
$ws.Range("I22").Value = 2983.3333
$ws.Range("J22").Value = 41676452
$ws.Range("K22").Value = 2983.3333
$ws.Range("L22").Value = 41676452
$ws.Range("M22").Value = -2688.3333
$ws.Range("N22").Value = -41677042

$ws.Range("H27").Value = 27785296
$ws.Range("I27").Value = 2983.3333
$ws.Range("J27").Value = 41676452
$ws.Range("K27").Value = 2983.3333
$ws.Range("L27").Value = 41676452
$ws.Range("M27").Value = -2876.3333
$ws.Range("N27").Value = -41676666

$ws.Range("H126").Value = 1611.2778
$ws.Range("I126").Value = 1437.125
$ws.Range("J126").Value = 3004.5
$ws.Range("K126").Value = 4311.375
$ws.Range("L126").Value = 9013.5
$ws.Range("M126").Value = -1841.375
$ws.Range("N126").Value = -13953.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 6613.75
$ws.Range("I45").Value = 6590
$ws.Range("J45").Value = 6621.6665
$ws.Range("K45").Value = 6590
$ws.Range("L45").Value = 6621.6665
$ws.Range("M45").Value = -6099
$ws.Range("N45").Value = -7603.6665

$ws.Range("H126").Value = 23110770
$ws.Range("I126").Value = 30220352
$ws.Range("J126").Value = 4625
$ws.Range("K126").Value = 90661056
$ws.Range("L126").Value = 13875
$ws.Range("M126").Value = -90658586
$ws.Range("N126").Value = -18815
